$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed by Excel as a
# number (e.g. "1.00" -> 1, "0.170" -> 0.17); force them to stay text by
# applying a text number format before writing the value, matching the
# original inline-string (text) cell content exactly.
$ws.Range('D2').Value = '69.277.52'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.848.42'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.82'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.18'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('D7').Value = '3.855.59'
$ws.Range('E7').Value = '  +2.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.35'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.25'
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '4.485.32'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('D16').Value = '3.831.90'
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').Value = '69.407.11'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.67'
$ws.Range('E19').Value = '  +7.41%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.30'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '489.72'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.725'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000159'
$ws.Range('E24').Value = '  +5.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.60'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  -1.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.24'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.98'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.02'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.43'
$ws.Range('E33').Value = '  +2.60%  '
$ws.Range('D34').Value = '3.990.62'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('D35').Value = '3.787.96'
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.141'
$ws.Range('E38').Value = '  +5.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.93'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.321'
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '440.05'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.60'
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.00'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.42'
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '143.26'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').Value = '2.847.94'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0357'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.88'
$ws.Range('E51').Value = '  +10.97%  '
